$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell E8 currently holds "Good Morning" (shared string index 10); change it
# to a new string "Good Morning2", which is appended as a brand new shared
# string entry (index 22) per the target diff.
$ws.Range("E8").Value = "Good Morning2"
